$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.62"
$ws.Range("E2").Value = "'-1.05%"
$ws.Range("D3").Value = "'36.32"
$ws.Range("E3").Value = "'-3.67%"
$ws.Range("D4").Value = "'5.113"
$ws.Range("E4").Value = "'-0.19%"
$ws.Range("D5").Value = "'0.07725"
$ws.Range("E5").Value = "'-2.21%"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "'4.386"
$ws.Range("E6").Value = "'-0.55%"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'8.295"
$ws.Range("E7").Value = "'0.45%"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").Value = "'1.851"
$ws.Range("E8").Value = "'-2.74%"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").Value = "'2.943"
$ws.Range("E9").Value = "'1.70%"
$ws.Range("B10").Value = "MXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D10").Value = "'0.9223"
$ws.Range("E10").Value = "'-0.62%"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.1125"
$ws.Range("E11").Value = "'-8.15%"
$ws.Range("B12").Value = "WazirX"
$ws.Range("C12").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D12").Value = "'0.1863"
$ws.Range("E12").Value = "'-3.15%"
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").Value = "'0.08764"
$ws.Range("E13").Value = "'-3.63%"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.03298"
$ws.Range("E14").Value = "'-0.96%"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09537"
$ws.Range("E15").Value = "'-0.99%"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001384"
$ws.Range("E16").Value = "'0.35%"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.005963"
$ws.Range("E17").Value = "'3.69%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.375"
$ws.Range("E18").Value = "'-4.62%"
$ws.Range("D19").Value = "'0.3448"
$ws.Range("E19").Value = "'1.16%"
$ws.Range("D20").Value = "'6.346"
$ws.Range("E20").Value = "'19.63%"
$ws.Range("D21").Value = "'0.1292"
$ws.Range("E21").Value = "'1.05%"
$ws.Range("D22").Value = "'0.2316"
$ws.Range("E22").Value = "'-11.46%"
$ws.Range("D23").Value = "'0.04347"
$ws.Range("E23").Value = "'-0.53%"
$ws.Range("D24").Value = "'0.001202"
$ws.Range("E24").Value = "'-2.93%"
$ws.Range("D25").Value = "'0.004259"
$ws.Range("E25").Value = "'-1.02%"
$ws.Range("D26").Value = "'0.0001203"
$ws.Range("E26").Value = "'-1.39%"
$ws.Range("D27").Value = "'0.0002907"
$ws.Range("D39").Value = "'0.02107"
$ws.Range("E39").Value = "'-0.79%"
$ws.Range("D40").Value = "'0.04901"
$ws.Range("E40").Value = "'-5.26%"
$ws.Range("D41").Value = "'0.007586"
$ws.Range("E41").Value = "'0.06%"
$ws.Range("E42").Value = "'-0.75%"
$ws.Range("D43").Value = "'0.008562"
$ws.Range("E43").Value = "'-6.31%"
$ws.Range("D44").Value = "'0.002073"
$ws.Range("E44").Value = "'1.16%"
$ws.Range("D45").Value = "'0.008524"
$ws.Range("E45").Value = "'-1.01%"
$ws.Range("D46").Value = "'0.00006533"
$ws.Range("E46").Value = "'-2.34%"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.23%"
$ws.Range("D48").Value = "'0.003303"
$ws.Range("E48").Value = "'16.48%"
$ws.Range("D49").Value = "'0.001446"
$ws.Range("E49").Value = "'20.53%"
$ws.Range("D50").Value = "'0.00002105"
$ws.Range("E50").Value = "'0.23%"
$ws.Range("D51").Value = "'0.0002005"
$ws.Range("E51").Value = "'0.23%"

Write-Host "Applied 102 cell updates"
